$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 316.84616
$ws.Range("I53").Value = 456.6
$ws.Range("J53").Value = 229.5
$ws.Range("K53").Value = 456.6
$ws.Range("L53").Value = 229.5
$ws.Range("M53").Value = 180.4
$ws.Range("N53").Value = -1503.5
$ws.Range("H64").Value = 8305.888999999999
$ws.Range("I64").Value = 7629
$ws.Range("J64").Value = 8499.286
$ws.Range("K64").Value = 7629
$ws.Range("L64").Value = 8499.286
$ws.Range("M64").Value = -7381
$ws.Range("N64").Value = -8995.286
$ws.Range("H67").Value = 8305.888999999999
$ws.Range("I67").Value = 7629
$ws.Range("J67").Value = 8499.286
$ws.Range("K67").Value = 7629
$ws.Range("L67").Value = 8499.286
$ws.Range("M67").Value = -6771
$ws.Range("N67").Value = -10215.286
$ws.Range("H70").Value = 2366.5557
$ws.Range("I70").Value = 1850
$ws.Range("J70").Value = 3399.6667
$ws.Range("K70").Value = 5550
$ws.Range("L70").Value = 10199.0001
$ws.Range("M70").Value = -5280
$ws.Range("N70").Value = -10739.0001
$ws.Range("H73").Value = 2366.5557
$ws.Range("I73").Value = 1850
$ws.Range("J73").Value = 3399.6667
$ws.Range("K73").Value = 5550
$ws.Range("L73").Value = 10199.0001
$ws.Range("M73").Value = -4614
$ws.Range("N73").Value = -12071.0001
$ws.Range("H111").Value = 1075.8
$ws.Range("I111").Value = 1094.75
$ws.Range("J111").Value = 1000
$ws.Range("K111").Value = 3284.25
$ws.Range("L111").Value = 3000
$ws.Range("M111").Value = -217.25
$ws.Range("N111").Value = -9134
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4156.4287
$ws.Range("I45").Value = 3849.1667
$ws.Range("K45").Value = 3849.1667
$ws.Range("M45").Value = -3472.1667
$ws.Range("H61").Value = 5428.905
$ws.Range("I61").Value = 5001.294
$ws.Range("K61").Value = 5001.294
$ws.Range("M61").Value = -4789.294
$ws.Range("H102").Value = 6703.6
$ws.Range("I102").Value = 6703.6
$ws.Range("K102").Value = 6703.6
$ws.Range("M102").Value = -5081.6
$ws.Range("H122").Value = 1960.1666
$ws.Range("I122").Value = 1930.1875
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 5790.5625
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -3340.5625
$ws.Range("N122").Value = -11500
$ws.Range("H136").Value = 5428.905
$ws.Range("I136").Value = 5001.294
$ws.Range("K136").Value = 15003.882
$ws.Range("M136").Value = -12453.882
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4600.558
$ws.Range("I134").Value = 2167.9375
$ws.Range("J134").Value = 11677.272
$ws.Range("K134").Value = 6503.8125
$ws.Range("L134").Value = 35031.81600000001
$ws.Range("M134").Value = -3968.8125
$ws.Range("N134").Value = -40101.81600000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1677.0952
$ws.Range("I22").Value = 492.1111
$ws.Range("K22").Value = 492.1111
$ws.Range("M22").Value = -142.1111
$ws.Range("H58").Value = 16789.3
$ws.Range("I58").Value = 7270.4287
$ws.Range("K58").Value = 7270.4287
$ws.Range("M58").Value = -7067.4287
$ws.Range("H62").Value = 3125
$ws.Range("J62").Value = 3500
$ws.Range("L62").Value = 3500
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 3125
$ws.Range("J65").Value = 3500
$ws.Range("L65").Value = 17500
$ws.Range("N65").Value = -23740
$ws.Range("H120").Value = 70000
$ws.Range("J120").Value = 70000
$ws.Range("L120").Value = 70000
$ws.Range("N120").Value = -77258
$ws.Range("H136").Value = 16789.3
$ws.Range("I136").Value = 7270.4287
$ws.Range("K136").Value = 21811.2861
$ws.Range("M136").Value = -19261.2861
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6088595.5
$ws.Range("I4").Value = 10769551
$ws.Range("J4").Value = 3352.7
$ws.Range("K4").Value = 32308653
$ws.Range("L4").Value = 10058.1
$ws.Range("M4").Value = -32308541
$ws.Range("N4").Value = -10282.1
$ws.Range("H7").Value = 7.5
$ws.Range("I7").Value = 7.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 22.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 89.5
$ws.Range("N7").ClearContents()
$ws.Range("H34").Value = 418.83334
$ws.Range("J34").Value = 900
$ws.Range("L34").Value = 2700
$ws.Range("N34").Value = -2868
$ws.Range("H39").Value = 3583.6428
$ws.Range("J39").Value = 4164.25
$ws.Range("L39").Value = 12492.75
$ws.Range("N39").Value = -13080.75
$ws.Range("H55").Value = 49415796
$ws.Range("J55").Value = 5656.6665
$ws.Range("L55").Value = 16969.9995
$ws.Range("N55").Value = -17323.9995
$ws.Range("H131").Value = 2325.3115
$ws.Range("J131").Value = 2325.3115
$ws.Range("L131").Value = 6975.934499999999
$ws.Range("N131").Value = -17055.9345
$ws.Range("H137").Value = 8635.875
$ws.Range("J137").Value = 9837.833000000001
$ws.Range("L137").Value = 29513.499
$ws.Range("N137").Value = -39713.499
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H55").Value = 35000
$ws.Range("I55").Value = 35000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 35000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -34673
$ws.Range("N55").ClearContents()
$ws.Range("H102").Value = 1450.9375
$ws.Range("I102").Value = 1526.2
$ws.Range("K102").Value = 1526.2
$ws.Range("M102").Value = 95.79999999999995
$ws.Range("H107").Value = 923
$ws.Range("I107").Value = 469.63635
$ws.Range("J107").Value = 1546.375
$ws.Range("K107").Value = 469.63635
$ws.Range("L107").Value = 1546.375
$ws.Range("M107").Value = 1450.36365
$ws.Range("N107").Value = -5386.375
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H117").Value = 90000
$ws.Range("J117").Value = 90000
$ws.Range("L117").Value = 90000
$ws.Range("N117").Value = -96884
$ws.Range("H125").Value = 30197.334
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -54920
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2823.6316
$ws.Range("I40").Value = 1963.25
$ws.Range("K40").Value = 1963.25
$ws.Range("M40").Value = -1827.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4250
$ws.Range("I62").Value = 4250
$ws.Range("K62").Value = 4250
$ws.Range("M62").Value = -3626
$ws.Range("H65").Value = 4250
$ws.Range("I65").Value = 4250
$ws.Range("K65").Value = 21250
$ws.Range("M65").Value = -18130
